# Daily attendance processing - 2026-01-16 23:03:14
# Swap the order of names in the "Recorded By" column (G) from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact value is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$col = 7  # Column G = "Recorded By"

$changed = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed 'Recorded By' cells in column G."
